$wb = $excel.ActiveWorkbook

# =====================================================================
# Sheet: MotorClaim_Insurer
#   - A2 default/selected insurer changes from ICICI to Royal Sundaram
#   - A6 spelling correction: HFDC ERGO -> HDFC ERGO
#   - selection moves from A2 to F7
# =====================================================================
$wsMotor = $wb.Worksheets.Item("MotorClaim_Insurer")
$wsMotor.Range("A2").Value = "Royal Sundaram"
$wsMotor.Range("A6").Value = "HDFC ERGO"

# =====================================================================
# Sheet: SuperAdmin
#   - A1 header renamed to include the new client login row
#   - A6 now holds the new test-client email, old "client@client.com"
#     moves down into a freshly inserted row 7 (same row style as before)
#   - list validation on A2 grows to include the new row
#   - column A gets a bit wider to fit the longer header text
# =====================================================================
$wsSuperAdmin = $wb.Worksheets.Item("SuperAdmin")

# Insert a new row before the old blank trailer row (row 7) so the
# existing "client@client.com" entry can be preserved one row lower,
# matching the same A/B/C styling pattern used by row 4 / row 6.
$wsSuperAdmin.Rows.Item(7).Insert()

$wsSuperAdmin.Range("A1").Value = "SUPERADMIN AND CLIENT USERNAME (0,0)"
$wsSuperAdmin.Range("A6").Value = "testclient1@catgroup.uk"

$wsSuperAdmin.Range("A4:C4").Copy()
$wsSuperAdmin.Range("A7:C7").PasteSpecial(-4122)
$wsSuperAdmin.Range("A7").Value = "client@client.com"
$wsSuperAdmin.Range("B7").Value = ""
$wsSuperAdmin.Application.CutCopyMode = $false

$wsSuperAdmin.Columns.Item(1).ColumnWidth = 48.33203125

$dvSuperAdmin = $wsSuperAdmin.Range("A2").Validation
$dvSuperAdmin.Modify(3, 1, 1, "=`$A`$3:`$A`$7")

# =====================================================================
# Sheet: Customer
#   - A2 sample username swapped for a different seeded address
#   - selection moves from B2 to A2
# =====================================================================
$wsCustomer = $wb.Worksheets.Item("Customer")
$wsCustomer.Range("A2").Value = "sumbadlet@eay.jp"

# =====================================================================
# Sheet: Insurer_Log
#   - A3 loses its redundant "applyFill" styling, matching A4's plain
#     bordered style (fixes the stray/unused style slot)
#   - selection moves from A10 to D8, and this sheet is no longer the
#     active tab (SuperAdmin takes over, see below)
# =====================================================================
$wsInsurerLog = $wb.Worksheets.Item("Insurer_Log")
$wsInsurerLog.Range("A4").Copy()
$wsInsurerLog.Range("A3").PasteSpecial(-4122)
$wsInsurerLog.Application.CutCopyMode = $false

# =====================================================================
# Re-apply per-sheet selections (each sheet keeps its own remembered
# selection independent of which tab ends up active).
# =====================================================================
$wsMotor.Range("F7").Select()
$wsCustomer.Range("A2").Select()

$wsRepairer = $wb.Worksheets.Item("Repairer")
$wsRepairer.Range("C26").Select()

$wsInsurerLog.Range("D8").Select()

# SuperAdmin becomes the active tab - select/activate it last so it
# sticks as the workbook's active sheet on save.
$wsSuperAdmin.Range("A2").Select()
$wsSuperAdmin.Activate()
